# KHL referee stats refresh (2025-12-22 run)
# Updates the per-referee running totals on the "Главные" (main) and
# "Линейные" (linesmen) sheets for officials who worked an additional
# game since the previous snapshot, and bumps every row's as_of_utc
# timestamp to the new snapshot time.

$wb = $excel.ActiveWorkbook

# Sheet order in this workbook: 1=Глоссарий, 2=Главные, 3=Линейные
$wsMain   = $wb.Worksheets.Item(2)
$wsLinear = $wb.Worksheets.Item(3)

$newTimestamp = "2025-12-22 04:42:07"

# ---------------------------------------------------------------------
# "Главные" (Главные судьи) sheet - updated stat rows
# ---------------------------------------------------------------------

# Row 7: C7=23, D7=333, E7=148, F7=185, G7=14.48, H7=6.43, I7=8.039999999999999, J7=69, K7=70, V7=18, W7=6
$wsMain.Range("C7").Value = 23
$wsMain.Range("D7").Value = 333
$wsMain.Range("E7").Value = 148
$wsMain.Range("F7").Value = 185
$wsMain.Range("G7").Value = 14.48
$wsMain.Range("H7").Value = 6.43
$wsMain.Range("I7").Value = 8.039999999999999
$wsMain.Range("J7").Value = 69
$wsMain.Range("K7").Value = 70
$wsMain.Range("V7").Value = 18
$wsMain.Range("W7").Value = 6

# Row 8: C8=32, D8=631, E8=328, F8=303, G8=19.72, H8=10.25, I8=9.470000000000001, J8=129, K8=129
$wsMain.Range("C8").Value = 32
$wsMain.Range("D8").Value = 631
$wsMain.Range("E8").Value = 328
$wsMain.Range("F8").Value = 303
$wsMain.Range("G8").Value = 19.72
$wsMain.Range("H8").Value = 10.25
$wsMain.Range("I8").Value = 9.470000000000001
$wsMain.Range("J8").Value = 129
$wsMain.Range("K8").Value = 129

# Row 16: C16=33, D16=581, F16=293, G16=17.61, H16=8.73, I16=8.880000000000001, K16=114
$wsMain.Range("C16").Value = 33
$wsMain.Range("D16").Value = 581
$wsMain.Range("F16").Value = 293
$wsMain.Range("G16").Value = 17.61
$wsMain.Range("H16").Value = 8.73
$wsMain.Range("I16").Value = 8.880000000000001
$wsMain.Range("K16").Value = 114

# Row 18: C18=34, D18=557, F18=287, G18=16.38, H18=7.94, I18=8.44, K18=121
$wsMain.Range("C18").Value = 34
$wsMain.Range("D18").Value = 557
$wsMain.Range("F18").Value = 287
$wsMain.Range("G18").Value = 16.38
$wsMain.Range("H18").Value = 7.94
$wsMain.Range("I18").Value = 8.44
$wsMain.Range("K18").Value = 121

# Row 19: C19=27, D19=454, E19=224, F19=230, G19=16.81, H19=8.300000000000001, I19=8.52, J19=107, K19=100
$wsMain.Range("C19").Value = 27
$wsMain.Range("D19").Value = 454
$wsMain.Range("E19").Value = 224
$wsMain.Range("F19").Value = 230
$wsMain.Range("G19").Value = 16.81
$wsMain.Range("H19").Value = 8.300000000000001
$wsMain.Range("I19").Value = 8.52
$wsMain.Range("J19").Value = 107
$wsMain.Range("K19").Value = 100

# Row 23: C23=22, D23=297, E23=120, F23=177, G23=13.5, H23=5.45, I23=8.050000000000001, J23=55, K23=71, M23=3, V23=6, W23=8
$wsMain.Range("C23").Value = 22
$wsMain.Range("D23").Value = 297
$wsMain.Range("E23").Value = 120
$wsMain.Range("F23").Value = 177
$wsMain.Range("G23").Value = 13.5
$wsMain.Range("H23").Value = 5.45
$wsMain.Range("I23").Value = 8.050000000000001
$wsMain.Range("J23").Value = 55
$wsMain.Range("K23").Value = 71
$wsMain.Range("M23").Value = 3
$wsMain.Range("V23").Value = 6
$wsMain.Range("W23").Value = 8

# Row 25: C25=35, D25=551, E25=266, F25=285, G25=15.74, H25=7.6, I25=8.140000000000001, J25=128, K25=135, V25=16, W25=20
$wsMain.Range("C25").Value = 35
$wsMain.Range("D25").Value = 551
$wsMain.Range("E25").Value = 266
$wsMain.Range("F25").Value = 285
$wsMain.Range("G25").Value = 15.74
$wsMain.Range("H25").Value = 7.6
$wsMain.Range("I25").Value = 8.140000000000001
$wsMain.Range("J25").Value = 128
$wsMain.Range("K25").Value = 135
$wsMain.Range("V25").Value = 16
$wsMain.Range("W25").Value = 20

# Row 26: C26=16, D26=334, E26=161, F26=173, G26=20.88, H26=10.06, I26=10.81, J26=63, K26=64
$wsMain.Range("C26").Value = 16
$wsMain.Range("D26").Value = 334
$wsMain.Range("E26").Value = 161
$wsMain.Range("F26").Value = 173
$wsMain.Range("G26").Value = 20.88
$wsMain.Range("H26").Value = 10.06
$wsMain.Range("I26").Value = 10.81
$wsMain.Range("J26").Value = 63
$wsMain.Range("K26").Value = 64

# ---------------------------------------------------------------------
# "Линейные" (линейные судьи / linesmen) sheet - updated stat rows
# ---------------------------------------------------------------------

# Row 8: C8=31, D8=496, E8=207, F8=289, G8=16, H8=6.68, I8=9.32, J8=86, K8=112, M8=5, V8=8, W8=8
$wsLinear.Range("C8").Value = 31
$wsLinear.Range("D8").Value = 496
$wsLinear.Range("E8").Value = 207
$wsLinear.Range("F8").Value = 289
$wsLinear.Range("G8").Value = 16
$wsLinear.Range("H8").Value = 6.68
$wsLinear.Range("I8").Value = 9.32
$wsLinear.Range("J8").Value = 86
$wsLinear.Range("K8").Value = 112
$wsLinear.Range("M8").Value = 5
$wsLinear.Range("V8").Value = 8
$wsLinear.Range("W8").Value = 8

# Row 9: C9=31, D9=573, E9=246, F9=327, G9=18.48, H9=7.94, I9=10.55, J9=108, K9=136
$wsLinear.Range("C9").Value = 31
$wsLinear.Range("D9").Value = 573
$wsLinear.Range("E9").Value = 246
$wsLinear.Range("F9").Value = 327
$wsLinear.Range("G9").Value = 18.48
$wsLinear.Range("H9").Value = 7.94
$wsLinear.Range("I9").Value = 10.55
$wsLinear.Range("J9").Value = 108
$wsLinear.Range("K9").Value = 136

# Row 12: C12=31, D12=515, E12=244, F12=271, G12=16.61, H12=7.87, I12=8.74, J12=112, K12=123, M12=5, V12=20, W12=14
$wsLinear.Range("C12").Value = 31
$wsLinear.Range("D12").Value = 515
$wsLinear.Range("E12").Value = 244
$wsLinear.Range("F12").Value = 271
$wsLinear.Range("G12").Value = 16.61
$wsLinear.Range("H12").Value = 7.87
$wsLinear.Range("I12").Value = 8.74
$wsLinear.Range("J12").Value = 112
$wsLinear.Range("K12").Value = 123
$wsLinear.Range("M12").Value = 5
$wsLinear.Range("V12").Value = 20
$wsLinear.Range("W12").Value = 14

# Row 14: C14=32, D14=514, E14=261, F14=253, H14=8.16, I14=7.91, J14=128, K14=119
$wsLinear.Range("C14").Value = 32
$wsLinear.Range("D14").Value = 514
$wsLinear.Range("E14").Value = 261
$wsLinear.Range("F14").Value = 253
$wsLinear.Range("H14").Value = 8.16
$wsLinear.Range("I14").Value = 7.91
$wsLinear.Range("J14").Value = 128
$wsLinear.Range("K14").Value = 119

# Row 16: C16=33, D16=512, F16=274, G16=15.52, H16=7.21, I16=8.300000000000001, K16=117
$wsLinear.Range("C16").Value = 33
$wsLinear.Range("D16").Value = 512
$wsLinear.Range("F16").Value = 274
$wsLinear.Range("G16").Value = 15.52
$wsLinear.Range("H16").Value = 7.21
$wsLinear.Range("I16").Value = 8.300000000000001
$wsLinear.Range("K16").Value = 117

# Row 18: C18=36, D18=604, E18=283, F18=321, G18=16.78, H18=7.86, I18=8.92, J18=134, K18=138
$wsLinear.Range("C18").Value = 36
$wsLinear.Range("D18").Value = 604
$wsLinear.Range("E18").Value = 283
$wsLinear.Range("F18").Value = 321
$wsLinear.Range("G18").Value = 16.78
$wsLinear.Range("H18").Value = 7.86
$wsLinear.Range("I18").Value = 8.92
$wsLinear.Range("J18").Value = 134
$wsLinear.Range("K18").Value = 138

# Row 21: C21=35, D21=667, E21=281, F21=386, G21=19.06, H21=8.029999999999999, I21=11.03, J21=128, K21=158
$wsLinear.Range("C21").Value = 35
$wsLinear.Range("D21").Value = 667
$wsLinear.Range("E21").Value = 281
$wsLinear.Range("F21").Value = 386
$wsLinear.Range("G21").Value = 19.06
$wsLinear.Range("H21").Value = 8.029999999999999
$wsLinear.Range("I21").Value = 11.03
$wsLinear.Range("J21").Value = 128
$wsLinear.Range("K21").Value = 158

# Row 22: C22=26, D22=518, E22=270, F22=248, G22=19.92, H22=10.38, I22=9.539999999999999, J22=105, K22=114, V22=18, W22=32
$wsLinear.Range("C22").Value = 26
$wsLinear.Range("D22").Value = 518
$wsLinear.Range("E22").Value = 270
$wsLinear.Range("F22").Value = 248
$wsLinear.Range("G22").Value = 19.92
$wsLinear.Range("H22").Value = 10.38
$wsLinear.Range("I22").Value = 9.539999999999999
$wsLinear.Range("J22").Value = 105
$wsLinear.Range("K22").Value = 114
$wsLinear.Range("V22").Value = 18
$wsLinear.Range("W22").Value = 32

# ---------------------------------------------------------------------
# Bump the as_of_utc snapshot timestamp (column AA) for every referee
# row (2..26) on both stat sheets.
# ---------------------------------------------------------------------
for ($row = 2; $row -le 26; $row++) {
    $wsMain.Range("AA$row").Value = $newTimestamp
    $wsLinear.Range("AA$row").Value = $newTimestamp
}
